$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25: Personal Leave note added under Application column (C)
$ws.Range("C25").Value = "Personal Leave"

# Row 26: Application + Comments added
$ws.Range("C26").Value = "Hayaai"
$ws.Range("D26").Value = "Dashboard requirement db queries writing going on"

# Row 27: Application + Comments added.
# C27 already contained an (empty) cell carrying the old "s=9" style, so a
# plain value assignment would keep that style. Copy the format from a
# neighboring cell that already carries the desired default style (s=1)
# before writing the value, so the on-disk style index matches.
$ws.Range("D27").Value = "Suppor for hayaai app for remove unwanted datas in table and dashboard queries going on."
$ws.Range("D27").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Value = "Hayaai"

# Row 30: Application + Comments added
$ws.Range("C30").Value = "Hayaai and Mujistore"
$ws.Range("D30").Value = "Dashboard requirement db queries going on and tested Mujistore issues"

# Row 31: Application + Comments added
$ws.Range("C31").Value = "Mujistore"
$ws.Range("D31").Value = "Fixing Mujistore issues"

# Restore the cursor/selection to where the author left it.
$ws.Range("D31").Select()
